$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.7151258441532331
$wsP.Range("D2").Value = 0.7443100016536195
$wsP.Range("E2").Value = 0.3223909503759028
$wsP.Range("F2").Value = 0.8078318456422755

$wsP.Range("B3").Value = 0.7151258441532331
$wsP.Range("D3").Value = 0.9715219719395836
$wsP.Range("E3").Value = 0.5067410597938991
$wsP.Range("F3").Value = 0.1350894018918429

$wsP.Range("B4").Value = 0.7443100016536195
$wsP.Range("C4").Value = 0.9715219719395836
$wsP.Range("E4").Value = 0.5229438200122272
$wsP.Range("F4").Value = 0.2535418009996488

$wsP.Range("B5").Value = 0.3223909503759028
$wsP.Range("C5").Value = 0.5067410597938991
$wsP.Range("D5").Value = 0.5229438200122272
$wsP.Range("F5").Value = 0.08083871874502124

$wsP.Range("B6").Value = 0.8078318456422755
$wsP.Range("C6").Value = 0.1350894018918429
$wsP.Range("D6").Value = 0.2535418009996488
$wsP.Range("E6").Value = 0.08083871874502124

# --- Sheet: Estadisticos_DM ---
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = 0.3697205411379265
$wsD.Range("D2").Value = 0.330284015152658
$wsD.Range("E2").Value = 1.012317899114806
$wsD.Range("F2").Value = -0.2461713678837381

$wsD.Range("B3").Value = -0.3697205411379265
$wsD.Range("D3").Value = -0.03610787837850123
$wsD.Range("E3").Value = 0.6749477410825059
$wsD.Range("F3").Value = -1.551333474605084

$wsD.Range("B4").Value = -0.330284015152658
$wsD.Range("C4").Value = 0.03610787837850123
$wsD.Range("E4").Value = 0.6491775104341889
$wsD.Range("F4").Value = -1.172486885121617

$wsD.Range("B5").Value = -1.012317899114806
$wsD.Range("C5").Value = -0.6749477410825059
$wsD.Range("D5").Value = -0.6491775104341889
$wsD.Range("F5").Value = -1.829968784687759

$wsD.Range("B6").Value = 0.2461713678837381
$wsD.Range("C6").Value = 1.551333474605084
$wsD.Range("D6").Value = 1.172486885121617
$wsD.Range("E6").Value = 1.829968784687759

$wb.Save()
